# Update the ITEM_ID values (column A) of the "users" sheet with a fresh
# batch of generated hash-style identifiers, as part of adding
# CartPage/CartTest/testng.xml. This introduces 50 new unique shared
# strings and rewrites A2:A51 to reference them instead of the previous
# identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

$newItemIds = @(
    "64cda32", "6557686", "5df29c9", "6363e33", "d3c4ef7", "c054e5e", "f96a497", "db7cfd0", "786f3d7", "7633f2f",
    "050c98d", "3a47835", "c6e056b", "6af4d48", "78f58cb", "c11c197", "8b2fc9b", "13eeabc", "561f00d", "ab16b44",
    "4f1bb5f", "b2c317c", "88ff7ce", "81631aa", "51e0737", "6b9b71d", "48849e1", "9265672", "1234223", "7bcd4b3",
    "2a17c26", "233f76d", "bf8a9d9", "a112f10", "9ad27a3", "0f2670d", "8e494eb", "3e67820", "cd1beb3", "240d155",
    "e131504", "896814c", "118359f", "76cbc16", "a9bb6e5", "5c5f64f", "e48e2d5", "617054d", "9e6fca1", "184222d"
)

$startRow = 2

# A handful of the new identifiers are purely numeric-looking (or parse as
# scientific notation), e.g. "6557686" or "6363e33". Left alone, assigning
# those through .Value would make Excel coerce them into numbers instead of
# text, same as typing them into the grid by hand. Force those specific
# cells to Text format first so they round-trip as the literal strings that
# belong in this ITEM_ID column, matching the rest of the column.
$numberLikeRows = @(3, 5, 26, 28, 29, 30, 39)

for ($i = 0; $i -lt $newItemIds.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    if ($numberLikeRows -contains $row) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $newItemIds[$i]
}

$wb.Save()
